$d = $word.ActiveDocument

# 1) Title change: "From White to Rainbow" -> "Understanding the Simple Electroscope"
$d.Content.Find.Execute("From White to Rainbow", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Understanding the Simple Electroscope", 2) | Out-Null

# 2) Consolidate the "aluminum" paragraph runs (removes the proofErr spell-check markers
#    split around "aluminum" by replacing the whole sentence span with itself, which
#    causes the engine to re-flow it into a single run).
$oldAluminum = "The rod should have a hook-like structure at the bottom so that the aluminum leaves can hang without falling easily. The leaves are typically made of aluminum foil and are attached at the bottom of the rod inside a glass jar."
$newAluminum = "The rod should have a hook-like structure at the bottom so that the aluminum leaves can hang without falling easily. The leaves are typically made of aluminum foil and are attached at the bottom of the rod inside a glass jar."
$d.Content.Find.Execute($oldAluminum, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newAluminum, 2) | Out-Null

# 3) Shift the split point between the two runs of the battery paragraph:
#    move "the voltage of a battery. This can be done by connecting one terminal of the
#    battery to the metal rod of the electroscope and " from the end of the first run to
#    the start of the second run.

# 3a) Trim the phrase off the end of the first run.
$oldRun1 = "it is possible to use an electroscope to detect the voltage of a battery. This can be done by connecting one terminal of the battery to the metal rod of the electroscope and "
$newRun1 = "it is possible to use an electroscope to detect "
$d.Content.Find.Execute($oldRun1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newRun1, 2) | Out-Null

# 3b) Insert the phrase right before "the other terminal" (start of the second run),
#     using a collapsed Range so the insertion lands inside the preceding run and keeps
#     the second run (with its lastRenderedPageBreak) intact.
$movedPhrase = "the voltage of a battery. This can be done by connecting one terminal of the battery to the metal rod of the electroscope and "
$fullText = $d.Content.Text
$insertAt = $fullText.IndexOf("the other terminal to the ground")
$insertRange = $d.Range($insertAt, $insertAt + 1)
$insertRange.Collapse(1)
$insertRange.Text = $movedPhrase
